$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column E = Solar (rows 13-26), divide by 1000
$solarRows = 13..26
foreach ($r in $solarRows) {
    $cell = $ws.Cells.Item($r, 5)
    $v = $cell.Value2
    if ($v -ne 0) {
        $cell.Value2 = $v / 1000
    }
}

# Column G = Wind, row 12, divide by 1000
$ws.Cells.Item(12, 7).Value2 = $ws.Cells.Item(12, 7).Value2 / 1000

# The whole data body (B2:G26) shares one number format (#,##0); the
# underlying numFmt definition is being changed to #,##0.0, which affects
# every cell using it (not only the Solar/Wind cells whose raw values moved).
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# Chart updates
$chart = $ws.ChartObjects(1).Chart
$valAxis = $chart.Axes(2)
$valAxis.AxisTitle.Text = "Kilowatts (kW)"
$valAxis.TickLabels.NumberFormat = "#,##0"
